# Generate Report for handback
# Update the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) timestamps on row 3
# of the zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-08 07:40:51"
$wsZhCn.Range("G3").Value = "2016-01-08 07:41:40"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-08 07:41:01"
$wsDeDe.Range("G3").Value = "2016-01-08 07:41:58"
